$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 202, shifting existing rows 202:336 down to 203:337
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new data record
$ws.Cells.Item(202, 1).Value = 3
$ws.Cells.Item(202, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(202, 3).Value = "Coquimbo"
$ws.Cells.Item(202, 4).Value = "02/04/2022"
$ws.Cells.Item(202, 5).Value = 5
$ws.Cells.Item(202, 6).Value = "Fruta"
$ws.Cells.Item(202, 7).Value = 100108
$ws.Cells.Item(202, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(202, 9).Value = 100108002
$ws.Cells.Item(202, 10).Value = "Mango"
$ws.Cells.Item(202, 11).Value = "Sin especificar"
$ws.Cells.Item(202, 12).Value = "Primera"
$ws.Cells.Item(202, 13).Value = 456
$ws.Cells.Item(202, 14).Value = 6500
$ws.Cells.Item(202, 15).Value = 7000
$ws.Cells.Item(202, 16).Value = 6750
$ws.Cells.Item(202, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(202, 18).Value = "Perú"
$ws.Cells.Item(202, 19).Value = 1688
$ws.Cells.Item(202, 20).Value = 4
